$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("router_vecinos")

# Insert a new "2 <-> 3" neighbor relationship pair into the router_vecinos table.
# The table lists neighbor pairs sorted by first column; inserting the pair (2,3)
# after the existing (2,*) rows and (3,2) after the existing (3,1) row keeps the
# ordering intact and shifts everything below down by one row each time.

# New row for (2,3), inserted after the last existing "A=2" row (old row 6)
$ws.Rows("6:6").Insert()
$ws.Range("A6:B6").ClearFormats()
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 3

# New row for (3,2), inserted after the existing (3,1) row (now at row 8)
$ws.Rows("9:9").Insert()
$ws.Range("A9:B9").ClearFormats()
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 2

# Make "router_vecinos" the active sheet/tab, with C9 as the selected cell
$ws.Activate() | Out-Null
$ws.Range("C9").Select() | Out-Null
